# Automatische test-sync: 2025-08-05 16:18:50
# Appends the newest logged test-mail ("Testmail #1: Ik heb nog geen geld terug.")
# as a new row (row 6) at the bottom of the Sheet1 log table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B6 needs an embedded line-break per paragraph, same convention Excel
# uses for in-cell line breaks (LF / Chr(10)).
$body = "Beste klant," + [char]10 + "Bedankt voor uw bericht. Ik begrijp dat u nog geen geld hebt ontvangen en ik help u hier graag mee verder. Om uw vraag goed te kunnen beantwoorden, zou ik wat meer informatie nodig hebben. Kunt u mij laten weten om welke transactie het gaat en eventueel het bijbehorende referentienummer? Op die manier kan ik het voor u nakijken en u verder helpen." + [char]10 + "Ik kijk uit naar uw reactie." + [char]10 + "Met vriendelijke groet," + [char]10 + "[Naam]" + [char]10 + "E-mailassistent"

$ws.Range("A6").Value = "Testmail #1: Ik heb nog geen geld terug."
$ws.Range("B6").Value = $body
$ws.Range("C6").Value = "Ik heb nog geen geld terug."
$ws.Range("D6").Value = "mailmind.test@zohomail.eu"
$ws.Range("E6").Value = "Retour / Terugbetaling"
$ws.Range("F6").Value = "2025-08-05 16:18:34"
$ws.Range("G6").Value = "Ja"
$ws.Range("H6").Value = "Nee"
$ws.Range("I6").Value = "Ja"
$ws.Range("J6").Value = "Nee"

# Keep the new row's height on the default (AutoFit recomputes it instead of
# leaving an inflated custom height from the multi-line text in B6).
$ws.Rows.Item(6).AutoFit()
